$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Section 1: rows 14-17 (Viewing list of vendors page test case) ---
$ws.Range("F14").Value = "It should not be displayed and should be as per parameter."
$ws.Range("E15").Value = "Validate companytID field."
$ws.Range("F15").Value = "It should not be displayed and should be as per parameter."

# --- Section 2: rows 37-40 (Adding new vendors form test case) ---
$ws.Range("F37").Value = "It should not be displayed and should be as per parameter."
$ws.Range("E38").Value = "Validate companytID field."
$ws.Range("F38").Value = "It should not be displayed and should be as per parameter."

# --- Section 3: rows 53-56 (Error validation test case) ---
$ws.Range("F53").Value = "It should not be displayed and should be as per parameter."
$ws.Range("E54").Value = "Validate companyID field."
$ws.Range("F54").Value = "It should not be displayed and should be as per parameter."

# Rows 53-54 lose their taller (29pt) row height in the edited workbook -
# reset back to the sheet's default (auto) row height.
$ws.Rows("53").AutoFit()
$ws.Rows("54").AutoFit()

# --- View state: scroll position + active selection ---
$ws.Range("F53:F54").Select()
$excel.ActiveWindow.ScrollRow = 37
